{"js": "// Remove the trailing \"Ver no Jupiter...\" and \"\u00a9 2020 ...\" footer\n// paragraphs (plus the blank paragraph separating them from the last\n// bibliography entry) that were scraped from the site chrome.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\n// Find the paragraph containing the \"Ver no Jupiter\" marker; the empty\n// paragraph immediately before it and the copyright paragraph right\n// after it are removed together with it.\nconst items = paragraphs.items;\nlet verIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === targetTexts[0]) {\n    verIndex = i;\n    break;\n  }\n}\n\nif (verIndex === -1) {\n  throw new Error(\"Could not locate the 'Ver no Jupiter' paragraph\");\n}\n\nconst copyrightIndex = verIndex + 1;\nif (items[copyrightIndex].text !== targetTexts[1]) {\n  throw new Error(\"Unexpected paragraph after 'Ver no Jupiter' marker\");\n}\n\nconst blankIndex = verIndex - 1;\nif (items[blankIndex].text !== \"\") {\n  throw new Error(\"Unexpected paragraph before 'Ver no Jupiter' marker\");\n}\n\n// Delete from the bottom up so earlier indices stay valid.\nitems[copyrightIndex].delete();\nitems[verIndex].delete();\nitems[blankIndex].delete();\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" and \"\u00a9 2020 ...\" footer\n# paragraphs (plus the blank paragraph separating them from the last\n# bibliography entry) that were scraped from the site chrome.\n$d = $word.ActiveDocument\n\n$verText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightText = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$paras = $d.Paragraphs\n$count = $paras.Count\n\n$verIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $paras.Item($i).Range.Text.TrimEnd([char]13)\n    if ($t -eq $verText) {\n        $verIndex = $i\n        break\n    }\n}\n\nif ($verIndex -eq -1) {\n    throw \"Could not locate the 'Ver no Jupiter' paragraph\"\n}\n\n$copyrightIndex = $verIndex + 1\n$copyrightActual = $paras.Item($copyrightIndex).Range.Text.TrimEnd([char]13)\nif ($copyrightActual -ne $copyrightText) {\n    throw \"Unexpected paragraph after 'Ver no Jupiter' marker\"\n}\n\n$blankIndex = $verIndex - 1\n$blankActual = $paras.Item($blankIndex).Range.Text.TrimEnd([char]13)\nif ($blankActual -ne \"\") {\n    throw \"Unexpected paragraph before 'Ver no Jupiter' marker\"\n}\n\n# Delete from the bottom up so earlier indices stay valid.\n$paras.Item($copyrightIndex).Range.Delete()\n$paras.Item($verIndex).Range.Delete()\n$paras.Item($blankIndex).Range.Delete()\n"}
